{"js": "// Dutch translation pass for \"WHAT IS SMARTCASH.docx\".\n// Each English sentence/run is located with a literal (match-case) search\n// and swapped for its Dutch counterpart via insertText(..., \"Replace\"),\n// which keeps the original run formatting (rPr) intact.\n\nconst NBSP = \"\\u00a0\";\n\nconst replacements = [\n  {\n    find: \"Payment after every 47500 blocks starting at 574100.\" + NBSP + \"Typically, around the 25th of each month.\",\n    replace: \"Betaling na 47500 blokken beginnend bij blok 574100.\" + NBSP + \"Normaal gesproken, rond de 25ste van elke maand.\"\n  },\n  {\n    find: \"All users need to move funds into addresses holding at least 1000 SMART before the snapshot to be counted.\",\n    replace: \"Alle gebruikers moeten geld verplaatsen naar adressen die ten minste 1000 SMART bevatten voordat de snapshot plaats vindt.\"\n  },\n  {\n    find: \"If you spend ANY amount from an address, it will be ineligible for SmartRewards until the next round.\",\n    replace: \"Als je WELK bedrag DAN OOK naar een adres overmaakt, komt het tot de volgende ronde niet meer in aanmerking voor SmartRewards.\"\n  },\n  {\n    find: \"InstantPay ensures SmartCash is always blazing fast, no matter what you are buying.\",\n    replace: \"InstantPay zorgt ervoor dat SmartCash razendsnel werkt, ongeacht wat je koopt.\"\n  },\n  {\n    find: NBSP + \"is going to be a key feature that allows instant\" + NBSP,\n    replace: NBSP + \"is een belangrijke functie waarmee je direct\" + NBSP\n  },\n  {\n    find: \"point-of-sale in-store and online purchases\",\n    replace: \"op locatie en online aankopen kunt doen\"\n  },\n  {\n    find: \". People should be able to use cryptocurrency as means to pay for their daily needs. InstantPay makes that a reality and allows for transactions to happen instantly, even faster than using your Visa or Mastercard.\",\n    replace: \" . Mensen moeten cryptogeld kunnen gebruiken om te betalen voor hun dagelijkse behoeften. InstantPay maakt dat een realiteit en zorgt ervoor dat transacties onmiddellijk kunnen plaatsvinden, zelfs sneller dan met je Visa of Mastercard.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n\n  // Replace every match (in this document each literal text is unique,\n  // but loop defensively in case of repeats).\n  for (const range of found.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Dutch translation pass for \"WHAT IS SMARTCASH.docx\".\n# Each English sentence/run is located with Find.Execute (exact case,\n# not whole-word since a couple of the spans start/end on a non-breaking\n# space rather than a word character) and swapped via ReplaceWith, which\n# keeps the original run formatting intact.\n\n$d = $word.ActiveDocument\n$nbsp = [char]0xa0\n\nfunction Replace-Text($range, [string]$findText, [string]$replaceText) {\n    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Find.Execute did not find: $findText\"\n    }\n}\n\n# 1) SmartRewards payout cadence sentence.\n$f1 = \"Payment after every 47500 blocks starting at 574100.\" + $nbsp + \"Typically, around the 25th of each month.\"\n$r1 = \"Betaling na 47500 blokken beginnend bij blok 574100.\" + $nbsp + \"Normaal gesproken, rond de 25ste van elke maand.\"\nReplace-Text $d.Content $f1 $r1\n\n# 2) Snapshot eligibility sentence.\n$f2 = \"All users need to move funds into addresses holding at least 1000 SMART before the snapshot to be counted.\"\n$r2 = \"Alle gebruikers moeten geld verplaatsen naar adressen die ten minste 1000 SMART bevatten voordat de snapshot plaats vindt.\"\nReplace-Text $d.Content $f2 $r2\n\n# 3) Spending-disqualifies-rewards sentence.\n$f3 = \"If you spend ANY amount from an address, it will be ineligible for SmartRewards until the next round.\"\n$r3 = \"Als je WELK bedrag DAN OOK naar een adres overmaakt, komt het tot de volgende ronde niet meer in aanmerking voor SmartRewards.\"\nReplace-Text $d.Content $f3 $r3\n\n# 4) InstantPay intro sentence.\n$f4 = \"InstantPay ensures SmartCash is always blazing fast, no matter what you are buying.\"\n$r4 = \"InstantPay zorgt ervoor dat SmartCash razendsnel werkt, ongeacht wat je koopt.\"\nReplace-Text $d.Content $f4 $r4\n\n# 5)-7) The \"InstantPay <i>...</i> purchases. People should...\" paragraph is\n# split across three runs; scope each Find to that paragraph so the other\n# (unchanged) \"InstantPay\" run at the start is left alone.\n$instantPayPara = $d.Paragraphs.Item(56)\n\n$f5 = $nbsp + \"is going to be a key feature that allows instant\" + $nbsp\n$r5 = $nbsp + \"is een belangrijke functie waarmee je direct\" + $nbsp\nReplace-Text $instantPayPara.Range $f5 $r5\n\n$f6 = \"point-of-sale in-store and online purchases\"\n$r6 = \"op locatie en online aankopen kunt doen\"\nReplace-Text $instantPayPara.Range $f6 $r6\n\n$f7 = \". People should be able to use cryptocurrency as means to pay for their daily needs. InstantPay makes that a reality and allows for transactions to happen instantly, even faster than using your Visa or Mastercard.\"\n$r7 = \" . Mensen moeten cryptogeld kunnen gebruiken om te betalen voor hun dagelijkse behoeften. InstantPay maakt dat een realiteit en zorgt ervoor dat transacties onmiddellijk kunnen plaatsvinden, zelfs sneller dan met je Visa of Mastercard.\"\nReplace-Text $instantPayPara.Range $f7 $r7\n"}
